# Remove column from alcohol data
# The measurement sheet (Sheet1) has an extra, redundant column (M) whose
# values duplicate the following column (N). Delete column M so that the
# old column N shifts left and becomes the new column M, trimming the
# sheet from A:N down to A:M.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the entire column M (13th column) - data to its right (old N)
# shifts left to take its place.
$ws.Columns.Item(13).Delete()

# Restore the active sheet/selection to reflect the new last column (M1),
# matching where the edit was made.
$ws.Activate()
$ws.Range("M1").Select() | Out-Null
